$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2
)

# 2. Prepend extra text to the "Programa resumido" paragraph
$d.Content.Find.Execute(
    "A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Estrutura e ultraestrutura dos materiais lignocelulósicos, celulose, hemiceluloses e outras polioses. Lignina, extrativos e composição da casca. Reações em meio ácido, meio alcalino. A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares.",
    2
)

# 3. Insert a space after the period before each numbered item (2-6) in the
#    Portuguese "Programa" paragraph only (leave the italic English paragraph
#    that follows untouched). Paragraph 14 holds this text.
$progPara = $d.Paragraphs.Item(14).Range
$progPara.Find.Execute("química.2.", $true, $false, $false, $false, $false, $true, 1, $false, "química. 2.", 2)
$progPara = $d.Paragraphs.Item(14).Range
$progPara.Find.Execute("papel.3.", $true, $false, $false, $false, $false, $true, 1, $false, "papel. 3.", 2)
$progPara = $d.Paragraphs.Item(14).Range
$progPara.Find.Execute("epóxidos.4.", $true, $false, $false, $false, $false, $true, 1, $false, "epóxidos. 4.", 2)
$progPara = $d.Paragraphs.Item(14).Range
$progPara.Find.Execute("vegetal.5.", $true, $false, $false, $false, $false, $true, 1, $false, "vegetal. 5.", 2)
$progPara = $d.Paragraphs.Item(14).Range
$progPara.Find.Execute("celular.6.", $true, $false, $false, $false, $false, $true, 1, $false, "celular. 6.", 2)

# 4. Update the final-grade formula
$d.Content.Find.Execute(
    "A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2 x 0,9 + Estudo de Caso x 0,1.",
    2
)

# 5. Update the recovery-grade formula
$d.Content.Find.Execute(
    "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2.",
    2
)
